# Corrected template for resident service.
# Inserts new "resident service" address-template rows (Province/City/Zone/
# Postal Code/Region) ahead of the gender/default rows, and appends the
# missing RPR_SUP_REJECT_EMAIL_SUBJECT translations for fra/ara/hin/kan/tam.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 5 new rows before row 1714 ---------------------------------
$ws.Rows.Item(1714).Insert()
$ws.Rows.Item(1714).Insert()
$ws.Rows.Item(1714).Insert()
$ws.Rows.Item(1714).Insert()
$ws.Rows.Item(1714).Insert()

function Set-Row($r, $lang, $code, $descr) {
    $ws.Range("A$r").Value = $lang
    $ws.Range("B$r").Value = $code
    $ws.Range("C$r").Value = $descr
    # D column carries the literal text "TRUE" (shared string), matching
    # the style/type used by every other row in the table - copy it from
    # the row just below (already populated, same s="1" text style) so the
    # inserted cell keeps the "is_active" text format instead of becoming
    # a native boolean.
    $ws.Range("D1719").Copy($ws.Range("D$r"))
}

Set-Row 1714 "eng" "mosip.province.template.property"    "Province"
Set-Row 1715 "eng" "mosip.city.template.property"         "City"
Set-Row 1716 "eng" "mosip.zone.template.property"         "Zone"
Set-Row 1717 "eng" "mosip.postal.code.template.property"  "Postal Code"
Set-Row 1718 "eng" "mosip.region.template.property"       "Region"

# --- Append 5 new rows after the existing data (now ends at row 1743) ---
function Set-NewRow($r, $lang, $code, $descr) {
    $ws.Range("A$r").Value = $lang
    $ws.Range("B$r").Value = $code
    $ws.Range("C$r").Value = $descr
    $ws.Range("D1743").Copy($ws.Range("D$r"))
}

Set-NewRow 1744 "fra" "RPR_SUP_REJECT_EMAIL_SUBJECT" "Template for Supervisor Reject Email Subject"
Set-NewRow 1745 "ara" "RPR_SUP_REJECT_EMAIL_SUBJECT" "Template for Supervisor Reject Email Subject"
Set-NewRow 1746 "hin" "RPR_SUP_REJECT_EMAIL_SUBJECT" "Template for Supervisor Reject Email Subject"
Set-NewRow 1747 "kan" "RPR_SUP_REJECT_EMAIL_SUBJECT" "Template for Supervisor Reject Email Subject"
Set-NewRow 1748 "tam" "RPR_SUP_REJECT_EMAIL_SUBJECT" "Template for Supervisor Reject Email Subject"

# --- Leave the selection where the author's last edit landed ----------
[void]$ws.Range("C1718").Select()
